# SVN Revision #7055 - Update Std design system maps and fan power per 3/3 CEC
# NACM system map document with corrections to fan power table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# TABLE T24N_2022BaseFanPwrIdx (rows 13-19): the system-type columns were
# re-mapped/expanded from 8 (D:K) to 9 (D:L) columns, splitting the old
# Sys1/Sys3/Sys7h/Sys7a groupings into the new Sys3a/Sys3b/Sys3c and
# Sys7a/Sys7b/Sys7c groupings, and moving Sys5/Sys6/Sys9 two columns to the
# right (J:L instead of G:I... see below for exact layout).
# ---------------------------------------------------------------------------

# Row 14: header / description labels
$ws.Range("D14").Value = "3a – SZAC"
$ws.Range("E14").Value = "3b – SZHP (no furnace)"
$ws.Range("F14").Value = "3c – SZDFHP (with furnace)"
$ws.Range("G14").Value = "7a – SZVAVAC "
$ws.Range("H14").Value = "7b – SZVAVHP"
$ws.Range("I14").Value = "7c – SZVAVDFHP (with furnace)"
$ws.Range("J14").Value = "5 – PVAV"
$ws.Range("K14").Value = "6 – VAV"
$ws.Range("L14").Value = "9 – HEATVENT"

# Row 14 used to have the wrap-text style (s="1") on every cell and an
# explicit ht="30"; now only A14 keeps that style and the row reverts to
# the default (un-customized) row height.
$ws.Range("D14:L14").Style = "Normal"
$ws.Rows(14).AutoFit()

# Row 15: system codes
$ws.Range("D15").Value = "Sys3a"
$ws.Range("E15").Value = "Sys3b"
$ws.Range("F15").Value = "Sys3c"
$ws.Range("G15").Value = "Sys7a"
$ws.Range("H15").Value = "Sys7b"
$ws.Range("I15").Value = "Sys7c"
$ws.Range("J15").Value = "Sys5"
$ws.Range("K15").Value = "Sys6"
$ws.Range("L15").Value = "Sys9"

# Row 16: FlowCap <=5000
$ws.Range("D16").Value = 0.802
$ws.Range("E16").Value = 0.744
$ws.Range("F16").Value = 0.802
$ws.Range("G16").Value = 0.802
$ws.Range("H16").Value = 0.744
$ws.Range("I16").Value = 0.802
$ws.Range("J16").Value = 1
$ws.Range("K16").Value = 0.977
$ws.Range("L16").Value = 0.616

# Row 17: FlowCap <=10000
$ws.Range("D17").Value = 0.78
$ws.Range("E17").Value = 0.72
$ws.Range("F17").Value = 0.78
$ws.Range("G17").Value = 0.78
$ws.Range("H17").Value = 0.72
$ws.Range("I17").Value = 0.78
$ws.Range("J17").Value = 1.022
$ws.Range("K17").Value = 1.013
$ws.Range("L17").Value = 0.62

# Row 18: FlowCap >10000
$ws.Range("D18").Value = 0.748
$ws.Range("E18").Value = 0.676
$ws.Range("F18").Value = 0.748
$ws.Range("G18").Value = 0.748
$ws.Range("H18").Value = 0.676
$ws.Range("I18").Value = 0.748
$ws.Range("J18").Value = 0.964
$ws.Range("K18").Value = 0.947
$ws.Range("L18").Value = 0.605

# ---------------------------------------------------------------------------
# Window / view cosmetics captured by the diff
# ---------------------------------------------------------------------------
$excel.ActiveWindow.Zoom = 70
$ws.Range("I20:I22").Select()
